$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.398.68'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.820.49'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.91'
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5108'
$ws.Range("E7").Value = '  -4.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3925'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07787'
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.73'
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.111'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.92'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.244'
$ws.Range("E14").Value = '  -1.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.470'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = '1.814.67'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001136'
$ws.Range("E17").Value = '  +5.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.47'
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06625'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.080'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").Value = '28.432.28'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.244'
$ws.Range("E25").Value = '  +3.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.17'
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").Value = '2.025.90'
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '155.03'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.401'
$ws.Range("E29").Value = '  -3.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.42'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("E31").Value = '  +0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.103'
$ws.Range("E32").Value = '  -2.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.658'
$ws.Range("E33").Value = '  -0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.646'
$ws.Range("E34").Value = '  -0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07041'
$ws.Range("E35").Value = '  -1.51%  '
$ws.Range("E36").Value = '  -2.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02319'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.183'
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.744'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6253'
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.17'
$ws.Range("E41").Value = '  -1.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.175'
$ws.Range("E42").Value = '  -1.27%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.384'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.41'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.727'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5878'
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.26'
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("E49").Value = '  -0.76%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("E51").Value = '  +0.00%  '
